$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (A8 / A9) ---
$ws.Range("A8").Value = "custom_redistr_satisfied"
$ws.Range("A9").Value = "custom_redistr_skip"

# --- Numeric re-run updates ---
$ws.Range("J4").Value = 4.42395634859102

$ws.Range("B5").Value = 242.045500729898
$ws.Range("C5").Value = 251.061082082337
$ws.Range("E5").Value = 238.646978806938
$ws.Range("H5").Value = 264.188486368283
$ws.Range("I5").Value = 246.596884457216
$ws.Range("J5").Value = 222.74051731422
$ws.Range("K5").Value = 232.973437908526
$ws.Range("M5").Value = 232.533807970671

$ws.Range("B6").Value = 5.11340902370102
$ws.Range("C6").Value = 5.39220439864498
$ws.Range("E6").Value = 5.10123556039544
$ws.Range("H6").Value = 5.65424149834322
$ws.Range("I6").Value = 5.27280465978773
$ws.Range("J6").Value = 4.68126628703672
$ws.Range("K6").Value = 4.6656688249704
$ws.Range("M6").Value = 4.90243297521537

$ws.Range("B7").Value = 46.3725250660407
$ws.Range("C7").Value = 41.1467443227155
$ws.Range("E7").Value = 46.8559069700612
$ws.Range("H7").Value = 35.8269923725672
$ws.Range("I7").Value = 54.8866021415646
$ws.Range("J7").Value = 73.9170977262703
$ws.Range("K7").Value = 36.9039796394692
$ws.Range("M7").Value = 56.5535147336437

$ws.Range("G8").Value = 56.6838694309219

$ws.Range("F9").Value = 34.7683813345703
